# Apply updated cryptocurrency price/volume data to Sheet1.
# Values that look like plain decimal numbers are written with a leading
# apostrophe (PowerShell: two single-quotes '' inside a single-quoted string)
# so Excel stores them as text (matching the source workbook, where these
# cells are text, e.g. "318.51", not the numeric value 318.51).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '48.031.99'
$ws.Range("E2").Value = '  +0.58%  '

$ws.Range("D3").Value = '2.495.51'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''318.51'
$ws.Range("E5").Value = '  -1.33%  '

$ws.Range("D6").Value = '''105.15'
$ws.Range("E6").Value = '  -2.54%  '

$ws.Range("D7").Value = '''0.521'
$ws.Range("E7").Value = '  -0.43%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = '''0.535'
$ws.Range("E9").Value = '  -3.94%  '

$ws.Range("D10").Value = '''38.69'
$ws.Range("E10").Value = '  -3.92%  '

$ws.Range("D11").Value = '''19.84'
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").Value = '''0.0801'
$ws.Range("E12").Value = '  -1.53%  '

$ws.Range("E13").Value = '  -0.80%  '

$ws.Range("D14").Value = '''7.02'
$ws.Range("E14").Value = '  -2.00%  '

$ws.Range("D15").Value = '2.888.28'
$ws.Range("E15").Value = '  -0.40%  '

$ws.Range("D16").Value = '2.501.96'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = '''0.828'
$ws.Range("E17").Value = '  -2.50%  '

$ws.Range("D18").Value = '47.904.34'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").Value = '''12.88'
$ws.Range("E19").Value = '  -3.33%  '

$ws.Range("E20").Value = '  +6.93%  '

$ws.Range("D21").Value = '''6.61'
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").Value = '0.0₃0931'
$ws.Range("E22").Value = '  -1.05%  '

$ws.Range("D23").Value = '''70.86'
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").Value = '''271.56'
$ws.Range("E24").Value = '  +9.85%  '

$ws.Range("D25").Value = '''2.51'
$ws.Range("E25").Value = '  -2.98%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").Value = '''25.62'
$ws.Range("E27").Value = '  -0.54%  '

$ws.Range("D28").Value = '''2.20'
$ws.Range("E28").Value = '  -3.20%  '

$ws.Range("D29").Value = '''0.142'
$ws.Range("E29").Value = '  +0.20%  '

$ws.Range("D30").Value = '''9.68'

$ws.Range("D31").Value = '''34.47'
$ws.Range("E31").Value = '  -1.20%  '

$ws.Range("D32").Value = '''49.27'

$ws.Range("E33").Value = '  -0.11%  '

$ws.Range("D34").Value = '''19.02'
$ws.Range("E34").Value = '  -4.96%  '

$ws.Range("D35").Value = '''5.25'
$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("D36").Value = '''0.0774'
$ws.Range("E36").Value = '  -1.21%  '

$ws.Range("D37").Value = '''1.93'
$ws.Range("E37").Value = '  -1.88%  '

$ws.Range("D38").Value = '''4.57'
$ws.Range("E38").Value = '  -2.28%  '

$ws.Range("E39").Value = '  -3.78%  '

$ws.Range("D40").Value = '''0.110'
$ws.Range("E40").Value = '  -1.47%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''22.21'
$ws.Range("E41").Value = '  -0.33%  '

$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '''2.21'
$ws.Range("E42").Value = '  +0.95%  '

$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''119.65'
$ws.Range("E43").Value = '  +0.96%  '

$ws.Range("D44").Value = '''0.0302'
$ws.Range("E44").Value = '  +1.56%  '

$ws.Range("D45").Value = '1.996.14'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").Value = '''3.18'
$ws.Range("E46").Value = '  +2.79%  '

$ws.Range("E47").Value = '  +4.23%  '

$ws.Range("D49").Value = '''8.88'
$ws.Range("E49").Value = '  -2.25%  '

$ws.Range("D50").Value = '''5.15'
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("D51").Value = '''78.39'
$ws.Range("E51").Value = '  +1.35%  '
